$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 423, shifting existing rows 423-513 down to 424-514.
$ws.Rows("423:423").Insert()

# Populate the newly inserted row 423 with its data (matches the formatting/values
# pattern used by every other data row in the sheet).
$ws.Range("A423").Value = 8
$ws.Range("B423").Value = "Terminal La Palmera de La Serena"
$ws.Range("C423").Value = "Coquimbo"
$ws.Range("D423").Value = 45275
$ws.Range("E423").Value = 4
$ws.Range("F423").Value = 100112012
$ws.Range("G423").Value = "Espinaca"
$ws.Range("H423").Value = "Sin especificar"
$ws.Range("I423").Value = "Primera"
$ws.Range("J423").Value = 1300
$ws.Range("K423").Value = 450
$ws.Range("L423").Value = 500
$ws.Range("M423").Value = 475
$ws.Range("N423").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O423").Value = "Provincia del Elquí"
$ws.Range("P423").Value = 950
$ws.Range("Q423").Value = 0.5
$ws.Range("R423").Value = "Hortaliza"
